# Update "想去人数" (want-to-go count) F-column figures to match the latest
# bilibili-sourced snapshot, flip one sheet2 ticket to "不可售", and add the
# newly-announced Jujutsu Kaisen cafe event on sheet3 (shifting the existing
# MyGO entry down a row with its own refreshed counter).
$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("展览"): F column (想去人数) refresh ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value2 = 212
$ws1.Cells.Item(5, 6).Value2 = 1730
$ws1.Cells.Item(6, 6).Value2 = 730
$ws1.Cells.Item(7, 6).Value2 = 2783
$ws1.Cells.Item(8, 6).Value2 = 2141
$ws1.Cells.Item(9, 6).Value2 = 891
$ws1.Cells.Item(10, 6).Value2 = 2367
$ws1.Cells.Item(11, 6).Value2 = 743
$ws1.Cells.Item(12, 6).Value2 = 6909
$ws1.Cells.Item(13, 6).Value2 = 144
$ws1.Cells.Item(15, 6).Value2 = 162
$ws1.Cells.Item(16, 6).Value2 = 1576
$ws1.Cells.Item(17, 6).Value2 = 1366
$ws1.Cells.Item(18, 6).Value2 = 1234
$ws1.Cells.Item(20, 6).Value2 = 2885
$ws1.Cells.Item(21, 6).Value2 = 2859
$ws1.Cells.Item(22, 6).Value2 = 2860
$ws1.Cells.Item(23, 6).Value2 = 830
$ws1.Cells.Item(24, 6).Value2 = 1139
$ws1.Cells.Item(25, 6).Value2 = 280
$ws1.Cells.Item(26, 6).Value2 = 5507
$ws1.Cells.Item(27, 6).Value2 = 305
$ws1.Cells.Item(29, 6).Value2 = 3837
$ws1.Cells.Item(30, 6).Value2 = 179
$ws1.Cells.Item(31, 6).Value2 = 647
$ws1.Cells.Item(32, 6).Value2 = 1742
$ws1.Cells.Item(33, 6).Value2 = 1094
$ws1.Cells.Item(34, 6).Value2 = 201
$ws1.Cells.Item(36, 6).Value2 = 93
$ws1.Cells.Item(37, 6).Value2 = 291
$ws1.Cells.Item(38, 6).Value2 = 1072
$ws1.Cells.Item(39, 6).Value2 = 441
$ws1.Cells.Item(40, 6).Value2 = 1784
$ws1.Cells.Item(41, 6).Value2 = 58
$ws1.Cells.Item(42, 6).Value2 = 268
$ws1.Cells.Item(44, 6).Value2 = 934
$ws1.Cells.Item(46, 6).Value2 = 528
$ws1.Cells.Item(47, 6).Value2 = 51
$ws1.Cells.Item(49, 6).Value2 = 65
$ws1.Cells.Item(50, 6).Value2 = 103

# --- Sheet 2 ("演出"): F column (想去人数) refresh ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(8, 6).Value2 = 504
$ws2.Cells.Item(12, 6).Value2 = 147
$ws2.Cells.Item(14, 6).Value2 = 968
$ws2.Cells.Item(18, 6).Value2 = 7
$ws2.Cells.Item(20, 6).Value2 = 620
$ws2.Cells.Item(21, 6).Value2 = 287
$ws2.Cells.Item(22, 6).Value2 = 365
$ws2.Cells.Item(28, 6).Value2 = 318
$ws2.Cells.Item(29, 6).Value2 = 92
$ws2.Cells.Item(33, 6).Value2 = 55
$ws2.Cells.Item(37, 6).Value2 = 227
$ws2.Cells.Item(41, 6).Value2 = 8

# Row 7 ticket is no longer purchasable - G column (最低票价) becomes text
$ws2.Cells.Item(7, 7).Value2 = "不可售"

# --- Sheet 3 ("本地生活"): F column (想去人数) refresh ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(4, 6).Value2 = 3310
$ws3.Cells.Item(7, 6).Value2 = 1491
$ws3.Cells.Item(9, 6).Value2 = 421
$ws3.Cells.Item(10, 6).Value2 = 2901
$ws3.Cells.Item(11, 6).Value2 = 347
$ws3.Cells.Item(12, 6).Value2 = 653
$ws3.Cells.Item(13, 6).Value2 = 791

# Insert the new "呪術廻戦カフェ2024第二季" row at 14, pushing the existing
# "BanG Dream! It's MyGO!!!!! x animate cafe" row down to 15.
$ws3.Rows.Item(14).Insert()

# Copy formatting (bold/bordered index style) from the row above into the new A14
$ws3.Cells.Item(13, 1).Copy($ws3.Cells.Item(14, 1))

$ws3.Cells.Item(14, 1).Value2 = 13
$ws3.Cells.Item(14, 2).Value2 = "'2024-07-27"
$ws3.Cells.Item(14, 3).Value2 = "上海·［呪術廻戦カフェ2024第二季］主题咖啡厅"
$ws3.Cells.Item(14, 4).Value2 = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$ws3.Cells.Item(14, 5).Value2 = "2024.07.27 00:00-08.31 23:59"
$ws3.Cells.Item(14, 6).Value2 = 263
$ws3.Cells.Item(14, 7).Value2 = 30
$ws3.Cells.Item(14, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89361"
$ws3.Cells.Item(14, 9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/cPSEorSS1721121089976.png"

# The shifted-down MyGO row (now row 15) gets its index bumped and a refreshed counter
$ws3.Cells.Item(15, 1).Value2 = 14
$ws3.Cells.Item(15, 6).Value2 = 1301

# --- Sheet 4 ("全部类型"): F column (想去人数) refresh ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value2 = 1491
$ws4.Cells.Item(4, 6).Value2 = 212
$ws4.Cells.Item(6, 6).Value2 = 421
$ws4.Cells.Item(7, 6).Value2 = 2901
$ws4.Cells.Item(8, 6).Value2 = 730
$ws4.Cells.Item(9, 6).Value2 = 2783
$ws4.Cells.Item(10, 6).Value2 = 347
$ws4.Cells.Item(11, 6).Value2 = 2141
$ws4.Cells.Item(12, 6).Value2 = 891
$ws4.Cells.Item(13, 6).Value2 = 2367
$ws4.Cells.Item(14, 6).Value2 = 147
$ws4.Cells.Item(15, 6).Value2 = 743
$ws4.Cells.Item(16, 6).Value2 = 6909
$ws4.Cells.Item(17, 6).Value2 = 144
$ws4.Cells.Item(18, 6).Value2 = 653
$ws4.Cells.Item(19, 6).Value2 = 791
$ws4.Cells.Item(20, 6).Value2 = 1576
$ws4.Cells.Item(21, 6).Value2 = 1366
$ws4.Cells.Item(22, 6).Value2 = 1234
$ws4.Cells.Item(24, 6).Value2 = 1301
$ws4.Cells.Item(25, 6).Value2 = 2885
$ws4.Cells.Item(26, 6).Value2 = 2860
$ws4.Cells.Item(27, 6).Value2 = 365
$ws4.Cells.Item(28, 6).Value2 = 830
$ws4.Cells.Item(29, 6).Value2 = 1139
$ws4.Cells.Item(30, 6).Value2 = 280
$ws4.Cells.Item(31, 6).Value2 = 5508
$ws4.Cells.Item(32, 6).Value2 = 305
$ws4.Cells.Item(33, 6).Value2 = 647
$ws4.Cells.Item(34, 6).Value2 = 318
$ws4.Cells.Item(35, 6).Value2 = 1742
$ws4.Cells.Item(36, 6).Value2 = 1094
$ws4.Cells.Item(38, 6).Value2 = 92
$ws4.Cells.Item(39, 6).Value2 = 93
$ws4.Cells.Item(40, 6).Value2 = 291
$ws4.Cells.Item(41, 6).Value2 = 1072
$ws4.Cells.Item(42, 6).Value2 = 441
$ws4.Cells.Item(43, 6).Value2 = 1784
$ws4.Cells.Item(44, 6).Value2 = 58
$ws4.Cells.Item(45, 6).Value2 = 268
$ws4.Cells.Item(47, 6).Value2 = 934
$ws4.Cells.Item(48, 6).Value2 = 528
$ws4.Cells.Item(49, 6).Value2 = 227
$ws4.Cells.Item(50, 6).Value2 = 227
$ws4.Cells.Item(51, 6).Value2 = 103

